$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.1800000000005
$ws.Range("G2").Value = 0.00181169300196804
$ws.Range("H2").Value = 0.01736370865652161
$ws.Range("K2").Value = 4.426578929815001
$ws.Range("L2").Value = "[1.5468907701493038, 7.3062670894806985]"
$ws.Range("M2").Value = 0.00267725356154136
$ws.Range("N2").Value = 0.005354507123082719
$ws.Range("O2").Value = -1.798789787641002
$ws.Range("P2").Value = "[-2.616421509296004, -0.9811580659860009]"
$ws.Range("Q2").Value = 0.00001944869760639456
$ws.Range("R2").Value = 0.00003889739521278912
$ws.Range("S2").Value = 13.2779275959799
$ws.Range("T2").Value = "[11.57444867123206, 14.981406520727733]"
$ws.Range("W2").Value = 7.208688688688831
$ws.Range("X2").Value = 3.932012012012088
$ws.Range("Y2").Value = 10.48536536536557

# Row 3 updates
$ws.Range("E3").Value = 24.32000000000036
$ws.Range("G3").Value = 0.02282489260903875
$ws.Range("H3").Value = 0.06381024366418274
$ws.Range("K3").Value = 3.398310193724882
$ws.Range("L3").Value = "[0.18949326739847194, 6.607127120051292]"
$ws.Range("M3").Value = 0.03799947492766553
$ws.Range("N3").Value = 0.03799947492766553
$ws.Range("O3").Value = 2.584974135386196
$ws.Range("P3").Value = "[1.1383949355350405, 4.031553335237352]"
$ws.Range("Q3").Value = 0.0005071179575535378
$ws.Range("R3").Value = 0.0005071179575535378
$ws.Range("S3").Value = 12.17181275145809
$ws.Range("T3").Value = "[10.452032014247159, 13.891593488669013]"
$ws.Range("W3").Value = 14.31447447447469
$ws.Range("X3").Value = 8.715275275275406
$ws.Range("Y3").Value = 19.91367367367397
